{"js": "// \"Add graph to forest\" \u2014 append a new bullet to the existing \"Notes\" list,\n// right after the paragraph ending \"... of 0.7 produced a reasonable fit.\",\n// describing the fit obtained by letting each tree predict on the full\n// dataset.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text.indexOf(\"of 0.7 produced a reasonable fit.\") !== -1) {\n    anchor = p;\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Could not find anchor paragraph for new list item.\");\n}\n\n// Inserting directly after an existing list item paragraph makes the new\n// paragraph continue the same bulleted list (same style + numPr) \u2014 exactly\n// like the other \"Running the random forest ...\" bullets above it.\nconst newPara = anchor.insertParagraph(\n  \"Allowing the each tree to predict on the complete dataset, not just the sample, produced a very good fit.\",\n  Word.InsertLocation.after\n);\n\nawait context.sync();\n", "ps1": "# \"Add graph to forest\" \u2014 append a new bullet to the existing \"Notes\" list,\n# right after the paragraph ending \"... of 0.7 produced a reasonable fit.\",\n# describing the fit obtained by letting each tree predict on the full\n# dataset.\n\n$d = $word.ActiveDocument\n\n# Locate the anchor text with Find, then resolve it to the containing\n# paragraph object (Find narrows $findRange to just the matched text, so we\n# walk $d.Paragraphs to find the paragraph that contains that match).\n$findRange = $d.Content\n$found = $findRange.Find.Execute(\"of 0.7 produced a reasonable fit.\")\n\nif (-not $found) {\n    throw \"Could not find anchor text for new list item.\"\n}\n\n$anchorStart = $findRange.Start\n$anchorEnd = $findRange.End\n\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Start -le $anchorStart -and $p.Range.End -ge $anchorEnd) {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -lt 0) {\n    throw \"Could not resolve anchor paragraph.\"\n}\n\n$target = $d.Paragraphs.Item($targetIndex)\n\n# Inserting a new paragraph mark right after this list item continues the\n# same bulleted list (same \"List Paragraph\" style + numPr) as the other\n# \"Running the random forest ...\" bullets above it.\n$target.Range.InsertParagraphAfter()\n\n$newPara = $d.Paragraphs.Item($targetIndex + 1)\n$newPara.Range.Text = \"Allowing the each tree to predict on the complete dataset, not just the sample, produced a very good fit.\"\n"}
